$wb = $excel.ActiveWorkbook

# "Template Setup" sheet - add a new data row (row 4) under the existing data.
$ws3 = $wb.Worksheets.Item("Template Setup")
$ws3.Range("A4").Value = 1567
$ws3.Range("A4").HorizontalAlignment = -4108   # xlCenter
$ws3.Range("A4").VerticalAlignment = -4160     # xlTop

# Move the active selection/tab to the "Template Setup" sheet, cell A2.
$ws3.Activate()
$ws3.Range("A2").Select()
